$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Cells.Item(2,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(2,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,6).Value = "Handed back: in sync with en-US"
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(3,3).Value = "Handed back: in sync with en-US"

$i2 = $zh.Cells.Item(2,9)
$i2.Value = "3bcb08c9-a9d8-4942-a370-17de35ffb070.md"
$i2.Font.Underline = 2
$i2.Font.Color = 15570276
$zh.Hyperlinks.Add($i2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95a1b82a7fa305a7a0e1dbe44d030c919cca0f8f/e2e/3bcb08c9-a9d8-4942-a370-17de35ffb070.md", "", "", "3bcb08c9-a9d8-4942-a370-17de35ffb070.md")
$zh.Cells.Item(2,10).Value = "3bcb08c9-a9d8-4942-a370-17de35ffb070.ec4b10587fecdf4c8eff71b14c94fc0a873e5364.zh-cn.xlf"
$zh.Cells.Item(2,11).Value = "2016-09-08 05:04:38"

$i3 = $zh.Cells.Item(3,9)
$i3.Value = "7910e7a4-dd47-4b90-a6a8-c79c1e20108d.md"
$i3.Font.Underline = 2
$i3.Font.Color = 15570276
$zh.Hyperlinks.Add($i3, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95a1b82a7fa305a7a0e1dbe44d030c919cca0f8f/e2e/7910e7a4-dd47-4b90-a6a8-c79c1e20108d.md", "", "", "7910e7a4-dd47-4b90-a6a8-c79c1e20108d.md")
$zh.Cells.Item(3,10).Value = "7910e7a4-dd47-4b90-a6a8-c79c1e20108d.4d0d0f75f7f64ceb0171ac88f3d1e43ba38e2ee0.zh-cn.xlf"
$zh.Cells.Item(3,11).Value = "2016-09-08 05:04:38"

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$de.Cells.Item(3,3).Value = "Handed back: in sync with en-US"

$j2 = $de.Cells.Item(2,9)
$j2.Value = "3bcb08c9-a9d8-4942-a370-17de35ffb070.md"
$j2.Font.Underline = 2
$j2.Font.Color = 15570276
$de.Hyperlinks.Add($j2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95a1b82a7fa305a7a0e1dbe44d030c919cca0f8f/e2e/3bcb08c9-a9d8-4942-a370-17de35ffb070.md", "", "", "3bcb08c9-a9d8-4942-a370-17de35ffb070.md")
$de.Cells.Item(2,10).Value = "3bcb08c9-a9d8-4942-a370-17de35ffb070.ec4b10587fecdf4c8eff71b14c94fc0a873e5364.de-de.xlf"
$de.Cells.Item(2,11).Value = "2016-09-08 05:04:47"

$j3 = $de.Cells.Item(3,9)
$j3.Value = "7910e7a4-dd47-4b90-a6a8-c79c1e20108d.md"
$j3.Font.Underline = 2
$j3.Font.Color = 15570276
$de.Hyperlinks.Add($j3, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95a1b82a7fa305a7a0e1dbe44d030c919cca0f8f/e2e/7910e7a4-dd47-4b90-a6a8-c79c1e20108d.md", "", "", "7910e7a4-dd47-4b90-a6a8-c79c1e20108d.md")
$de.Cells.Item(3,10).Value = "7910e7a4-dd47-4b90-a6a8-c79c1e20108d.4d0d0f75f7f64ceb0171ac88f3d1e43ba38e2ee0.de-de.xlf"
$de.Cells.Item(3,11).Value = "2016-09-08 05:04:47"

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "done"
